$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "P4TSSOP8_DEFROST"
$ws.Range("E14").Select()
